$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "07/30/2021"
$ws.Range("H2").Value = 'Española'
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("P2").Value = 600
$ws.Range("Q2").Value = 30

# Row 3
$ws.Range("D3").Value = "07/15/2022"
$ws.Range("H3").Value = 'Argentina(o)'
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 18000
$ws.Range("N3").Value = '$/caja 50 unidades'
$ws.Range("P3").Value = 360
$ws.Range("Q3").Value = 50

# Row 4
$ws.Range("D4").Value = "08/09/2022"
$ws.Range("J4").Value = 90

# Row 5
$ws.Range("D5").Value = "10/04/2022"
$ws.Range("H5").Value = 'Española'
$ws.Range("J5").Value = 150
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = '$/caja 30 unidades'
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 30

# Row 6
$ws.Range("D6").Value = "08/10/2021"
$ws.Range("H6").Value = 'Española'
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 16000
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("P6").Value = 533
$ws.Range("Q6").Value = 30

# Row 7
$ws.Range("D7").Value = "10/07/2022"
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11533
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 384
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = "08/05/2022"
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("J8").Value = 160
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 388
$ws.Range("Q8").Value = 40

# Row 9
$ws.Range("D9").Value = "09/07/2021"
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 16000
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("Q9").Value = 40

# Row 10
$ws.Range("D10").Value = "07/13/2021"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 16000
$ws.Range("P10").Value = 533

# Row 11
$ws.Range("D11").Value = "07/23/2021"
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 500

# Row 12
$ws.Range("D12").Value = "11/27/2020"
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("P12").Value = 250

# Row 13
$ws.Range("D13").Value = "08/13/2021"
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 16500
$ws.Range("L13").Value = 16500
$ws.Range("M13").Value = 16500
$ws.Range("P13").Value = 550

# Row 14
$ws.Range("D14").Value = "09/23/2022"
$ws.Range("H14").Value = 'Madrigal'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 12000
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("P14").Value = 300
$ws.Range("Q14").Value = 40

# Row 15
$ws.Range("D15").Value = "09/27/2022"
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("J15").Value = 180
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12444
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("P15").Value = 311
$ws.Range("Q15").Value = 40

# Row 16
$ws.Range("D16").Value = "11/26/2020"
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = 11000
$ws.Range("N16").Value = '$/caja 40 unidades'
$ws.Range("P16").Value = 275
$ws.Range("Q16").Value = 40

# Row 17
$ws.Range("D17").Value = "10/22/2021"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 11000
$ws.Range("M17").Value = 11000
$ws.Range("P17").Value = 275

# Row 18
$ws.Range("D18").Value = "07/09/2021"
$ws.Range("H18").Value = 'Española'
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 15000
$ws.Range("N18").Value = '$/caja 30 unidades'
$ws.Range("P18").Value = 500
$ws.Range("Q18").Value = 30

# Row 19
$ws.Range("D19").Value = "12/11/2020"
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 11000
$ws.Range("M19").Value = 11000
$ws.Range("P19").Value = 275

# Row 20
$ws.Range("D20").Value = "10/14/2022"
$ws.Range("J20").Value = 200

# Row 21
$ws.Range("D21").Value = "08/16/2022"
$ws.Range("J21").Value = 80

# Row 22
$ws.Range("D22").Value = "09/02/2022"
$ws.Range("I22").Value = 'Segunda'
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 13000
$ws.Range("N22").Value = '$/caja 50 unidades'
$ws.Range("P22").Value = 260
$ws.Range("Q22").Value = 50

# Row 23
$ws.Range("D23").Value = "07/29/2022"
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 16000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 16000
$ws.Range("N23").Value = '$/caja 40 unidades'
$ws.Range("P23").Value = 400
$ws.Range("Q23").Value = 40

# Row 24
$ws.Range("D24").Value = "10/12/2021"
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 120
$ws.Range("L24").Value = 11000
$ws.Range("M24").Value = 11000
$ws.Range("N24").Value = '$/caja 50 unidades'
$ws.Range("P24").Value = 220
$ws.Range("Q24").Value = 50

# Row 25
$ws.Range("D25").Value = "09/06/2022"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 13500
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13750
$ws.Range("P25").Value = 344

# Row 26
$ws.Range("D26").Value = "12/01/2020"
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 10000
$ws.Range("P26").Value = 250

# Row 27
$ws.Range("D27").Value = "10/19/2021"
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 12000
$ws.Range("N27").Value = '$/caja 40 unidades'
$ws.Range("P27").Value = 300
$ws.Range("Q27").Value = 40

# Row 28
$ws.Range("D28").Value = "08/30/2022"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14500
$ws.Range("P28").Value = 362

# Row 29
$ws.Range("D29").Value = "10/15/2021"
$ws.Range("J29").Value = 110
$ws.Range("K29").Value = 11000
$ws.Range("L29").Value = 11000
$ws.Range("M29").Value = 11000
$ws.Range("N29").Value = '$/caja 50 unidades'
$ws.Range("P29").Value = 220
$ws.Range("Q29").Value = 50

# Row 30
$ws.Range("D30").Value = "10/26/2021"
$ws.Range("H30").Value = 'Madrigal'
$ws.Range("J30").Value = 130
$ws.Range("K30").Value = 11000
$ws.Range("L30").Value = 11000
$ws.Range("M30").Value = 11000
$ws.Range("N30").Value = '$/caja 40 unidades'
$ws.Range("P30").Value = 275
$ws.Range("Q30").Value = 40

# Row 31
$ws.Range("D31").Value = "08/26/2022"
$ws.Range("H31").Value = 'Madrigal'
$ws.Range("J31").Value = 70
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("N31").Value = '$/caja 40 unidades'
$ws.Range("P31").Value = 375
$ws.Range("Q31").Value = 40

# Row 32
$ws.Range("D32").Value = "08/02/2022"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = 16500
$ws.Range("P32").Value = 412

# Row 33
$ws.Range("D33").Value = "09/28/2022"
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 12500
$ws.Range("P33").Value = 312

# Row 34
$ws.Range("D34").Value = "10/11/2022"
$ws.Range("H34").Value = 'Española'
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 10000
$ws.Range("N34").Value = '$/caja 30 unidades'
$ws.Range("P34").Value = 333
$ws.Range("Q34").Value = 30
